# reopening_ny.xlsx -- "Add files via upload" re-edit
#
# This re-applies the view-state + data changes described by the commit's
# canonical-XML diff:
#   - new_york_regions!G24:G30, G31:G32, G39:G43 gain a "phase_4_date"
#     value (2 weeks after the corresponding phase_3_date in column F),
#     formatted the same way as the existing date cells in column F.
#   - the active sheet/tab switches from new_york_phases (sheet 1) to
#     new_york_regions (sheet 2).
#   - each sheet's remembered scroll position / selection is updated to
#     match the author's final cursor position.
#
# Note: the x15ac:absPath (C:\ei\mobility\ -> C:\ei\replicate\) attribute
# is an MS-Office-internal, autosave-tracking artifact tied to the
# authoring machine's on-disk path. It is not exposed anywhere in the
# Excel object model (no Workbook/Application property round-trips into
# it), so it cannot be changed from automation code -- it is left as-is.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # new_york_phases
$ws2 = $wb.Worksheets.Item(2)   # new_york_regions

# ---------------------------------------------------------------------
# 1. Fill in the new "phase_4_date" (column G) values on new_york_regions.
#    Each new cell is two weeks (14 days) after its row's phase_3_date
#    (column F). Copy the F-column cell's format first so the new G cell
#    picks up the same date number-format style as its neighbours.
# ---------------------------------------------------------------------
$phase4Rows = 24..32 + 39..43

foreach ($r in $phase4Rows) {
    $src = $ws2.Range("F$r")
    $dst = $ws2.Range("G$r")

    $src.Copy()
    $dst.PasteSpecial(-4122)        # xlPasteFormats
    $dst.Value2 = $src.Value2 + 14  # Value2 keeps the raw date serial (avoids locale string round-tripping)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Update the remembered view/selection state for each sheet, then
#    leave new_york_regions as the active (selected) tab -- matching the
#    final cursor positions recorded in the workbook.
# ---------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("C29").Select()

$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("G43").Select()
